# Auto-generated edit script for tipo_financiador.xlsx
# Implements: append financiador/tipo-cobertura rows 1040-1069, highlight
# three rows in yellow, extend the AutoFilter / _FilterDatabase range, and
# reset the sheet view (scroll position + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ row=1040; a="Galeno ART - GALENO ASEGURADORA DE RIESGOS DEL TRABAJO SOCIEDAD ANONIMA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1041; a="TRINIDAD - SANATORIO DE LA TRINIDAD"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1042; a="''- SECRETARIA DE DESARROLLO SOCIAL"; b="Sin Cobertura"; astyle=5 },
    @{ row=1043; a="- SECRETARIA DE DESARROLLO SOCIAL"; b="Sin Cobertura"; astyle=5 },
    @{ row=1044; a="INCLUIR SALUD CORDOBA - INCLUIR SALUD CORDOBA"; b="Incluir Extra Cápita"; astyle=0 },
    @{ row=1045; a="DELEGACION METR - DELEGACION METROPOLITANA S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1046; a="INCLUIR SALUD CHUBUT - INCLUIR SALUD CHUBUT"; b="Incluir Extra Cápita"; astyle=0 },
    @{ row=1047; a="TIEMPOMEDICO - TIEMPO MEDICO S.R.L."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1048; a="''- UNIDAD DE GESTIO OPERATIVA FERROVIARIA DE EMERGENC"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1049; a="MEDICINA ESENCI - MEDICINA ESENCIAL S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1050; a="INCLUIR OTROS - INCLUIR SALUD OTROS"; b="Incluir Extra Cápita"; astyle=5 },
    @{ row=1051; a="AMPSI - ASOCIACION MUTUAL DE PSICOLOGOS AMPSI"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1052; a="PREMEDIC - GRUPO PREMEDIC S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1053; a="MINIST.SALUD SA - MINISTERIO DE SALUD PUBLICA DEL GOBIERNO DE LA PROVINCIA DE SALTA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1054; a="MERCANTIL ANDIN - COMPAÑÍA DE SEGUROS LA MERCANTIL ANDINA S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1055; a="CONSUL.BOLIVIA - CONSULADO GENERAL DE BOLIVIA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1056; a="OS PORTUARIOS P - OBRA SOCIAL PORTUARIOS DE PUERTO SAN MARTIN Y BELLA VISTA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1057; a="OBRA SOCIAL DE DIRECCION WITCE - OBRA SOCIAL DE DIRECCION WITCEL"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1058; a="INCHCAPE SHIPPI - INCHCAPE SHIPPING SERVICES ARGENTINA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1059; a="PREMEDIC - GRUPO PREMEDIC S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1060; a="PCIA SALUD HTAL - PROVINCIA SALUD HOSPITAL FRANCES"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1061; a="MERCANTIL ANDIN - COMPAÑÍA DE SEGUROS LA MERCANTIL ANDINA S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1062; a="PREMEDIC - GRUPO PREMEDIC S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1063; a="CONSUL.UCRANIA - CONSULADO GENERAL DE UCRANIA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1064; a="LUZ ART SA - LUZ ART SA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1065; a="OSPAISM - OBRA SOCIAL DEL PERSONAL DEL AZUCAR INGENIO SAN MARTIN"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1066; a="Galeno ART - GALENO ASEGURADORA DE RIESGOS DEL TRABAJO SOCIEDAD ANONIMA"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1067; a="''- UNIDAD DE GESTIO OPERATIVA FERROVIARIA DE EMERGENC"; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1068; a="MERCANTIL ANDIN - COMPAÑÍA DE SEGUROS LA MERCANTIL ANDINA S.A."; b="OOSS y Prepagas"; astyle=0 },
    @{ row=1069; a="MERCANTIL ANDIN - COMPAÑÍA DE SEGUROS LA MERCANTIL ANDINA S.A."; b="OOSS y Prepagas"; astyle=0 }
)

foreach ($item in $newRows) {
    $ws.Cells.Item($item.row, 1).Value2 = $item.a
    $ws.Cells.Item($item.row, 2).Value2 = $item.b
}

# Highlight the three "Sin Cobertura" / "Incluir Otros" rows in yellow,
# mirroring the workbook's existing conditional-formatting-derived dxf/fill
# so the styles table grows the same way it did in the authored edit.
$yellowRows = @(1042, 1043, 1050)
foreach ($r in $yellowRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.FormatConditions.Add(1, 1, "1") | Out-Null
    $cell.FormatConditions.Item(1).Interior.Color = 65535
    $cell.FormatConditions.Delete()
    $cell.Interior.Color = 65535
}

# Extend the AutoFilter range to cover the newly appended rows. Toggling the
# range's AutoFilter twice turns it off then back on over the new extent.
$lastRow = 1069
$filterRange = $ws.Range("A1:B" + $lastRow)
$filterRange.AutoFilter()
$filterRange.AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$fd = $wb.Names.Item("Hoja1!_FilterDatabase")
$fd.RefersTo = "=Hoja1!`$A`$1:`$B`$" + $lastRow

# Reset the view: scroll back to the top-left corner and select A7, matching
# the state the workbook was saved in after the edit.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A7").Select()
